# dokończona edycja i usuwanie inw
#
# Helper: force a value to be written as TEXT (inlineStr/shared-string in the
# saved file) without leaving a lingering custom number-format style behind.
function Set-TextCell {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value2 = $Text
    $Cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the old (now-invalid) inventory rows 17-20, but keep the date cell
# in column B (it keeps its number format / style, just loses its value).
# ---------------------------------------------------------------------------
$ws.Range("A17").ClearContents()
$ws.Range("C17:D17").ClearContents()
$ws.Range("B17").ClearContents()

$ws.Range("A18").ClearContents()
$ws.Range("C18:D18").ClearContents()
$ws.Range("B18").ClearContents()

$ws.Range("A19").ClearContents()
$ws.Range("C19:D19").ClearContents()
$ws.Range("B19").ClearContents()

$ws.Range("A20").ClearContents()
$ws.Range("C20:E20").ClearContents()
$ws.Range("B20").ClearContents()

# ---------------------------------------------------------------------------
# Add the new inventory rows 21-23.
# ---------------------------------------------------------------------------

# Row 21: piwo Żywiec
Set-TextCell $ws.Range("A21") "5901359000537"
$ws.Range("B21").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B21").Value2 = 45813
Set-TextCell $ws.Range("C21") "800"
Set-TextCell $ws.Range("D21") "3"
Set-TextCell $ws.Range("E21") "piwo Żywiec"

# Row 22: piwo Tyskie (barcode kept as a real number)
$ws.Range("A22").Value2 = 5901359072145
$ws.Range("B22").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B22").Value2 = 45813
Set-TextCell $ws.Range("C22") "600"
Set-TextCell $ws.Range("D22") "2"
Set-TextCell $ws.Range("E22") "piwo Tyskie"

# Row 23: piwo Lech Premium
Set-TextCell $ws.Range("A23") "5901359112568"
$ws.Range("B23").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B23").Value2 = 45813
Set-TextCell $ws.Range("C23") "800"
Set-TextCell $ws.Range("D23") "1"
Set-TextCell $ws.Range("E23") "piwo Lech Premium"

# ---------------------------------------------------------------------------
# Update the view: scroll so row 11 is the top visible row, and select A22
# (mirrors the cursor/scroll position left behind by the author's edit).
# ---------------------------------------------------------------------------
$null = $ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$null = $ws.Range("A22").Select()
